# Weekly refresh of the "Fruta, Vega Central Mapocho de Santiago - Breva"
# price sheet: each data row (2-13) is rewritten with that week's
# observation (date, quality, volume, prices, unit, origin, $/kg, kg/unit).
# Columns A,B,C,E,F,G,H,I,J,K (market/product identifiers) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  D = 44572; L = 'Primera';  M = 65;  N = 20000; O = 20000; P = 20000; Q = '$/bandeja 6 kilos'; R = 'Región Metropolitana';                 S = 3333; T = 6 }
    @{ Row = 3;  D = 44558; L = 'Especial'; M = 20;  N = 22000; O = 22000; P = 22000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 3667; T = 6 }
    @{ Row = 4;  D = 44558; L = 'Primera';  M = 25;  N = 18000; O = 18000; P = 18000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 3000; T = 6 }
    @{ Row = 5;  D = 44204; L = 'Primera';  M = 110; N = 7000;  O = 7500;  P = 7318;  Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 1045; T = 7 }
    @{ Row = 6;  D = 44550; L = 'Primera';  M = 60;  N = 24000; O = 24000; P = 24000; Q = '$/bandeja 7 kilos'; R = 'Región Metropolitana';                 S = 3429; T = 7 }
    @{ Row = 7;  D = 44561; L = 'Primera';  M = 200; N = 18000; O = 18000; P = 18000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 3000; T = 6 }
    @{ Row = 8;  D = 44189; L = 'Especial'; M = 20;  N = 15000; O = 15000; P = 15000; Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 2143; T = 7 }
    @{ Row = 9;  D = 44189; L = 'Primera';  M = 30;  N = 13000; O = 13000; P = 13000; Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 1857; T = 7 }
    @{ Row = 10; D = 44553; L = 'Especial'; M = 200; N = 22000; O = 22000; P = 22000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 3667; T = 6 }
    @{ Row = 11; D = 44553; L = 'Primera';  M = 150; N = 18000; O = 18000; P = 18000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 3000; T = 6 }
    @{ Row = 12; D = 44187; L = 'Especial'; M = 45;  N = 14000; O = 14000; P = 14000; Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 2000; T = 7 }
    @{ Row = 13; D = 44187; L = 'Primera';  M = 50;  N = 12000; O = 12000; P = 12000; Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua'; S = 1714; T = 7 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("N$n").Value = $r.N
    $ws.Range("O$n").Value = $r.O
    $ws.Range("P$n").Value = $r.P
    $ws.Range("Q$n").Value = $r.Q
    $ws.Range("R$n").Value = $r.R
    $ws.Range("S$n").Value = $r.S
    $ws.Range("T$n").Value = $r.T
}
